$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 28 ("o1-o1_1"/"cure_day7"/"o1_1"),
# which pushes it (and everything below) down by one row.
$ws.Rows.Item(28).Insert()

# New row 28: follow-up export entry for "o1-o1_2a" / "location_death_day7" / "o1_2a",
# highlighted in yellow for duplicate review.
$ws.Cells.Item(28, 1).Value = "o1-o1_2a"
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = "location_death_day7"
$ws.Cells.Item(28, 4).Value = "o1_2a"

$ws.Cells.Item(28, 1).Interior.Color = 65535
$ws.Cells.Item(28, 1).HorizontalAlignment = -4108
$ws.Cells.Item(28, 1).VerticalAlignment = -4108

$ws.Cells.Item(28, 3).Interior.Color = 65535
$ws.Cells.Item(28, 3).HorizontalAlignment = -4108
$ws.Cells.Item(28, 3).VerticalAlignment = -4108

# Update the view: scroll down a bit and move the active selection.
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("G20").Select()
